$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 15 (old "Hexaplex trunculus / HEXATRU" entry in the 2-RAP block);
# remaining rows shift up by one, so the former row 16 becomes row 15, etc.,
# and the former last row (47) disappears, shrinking the used range to A1:K46.
$ws.Cells.Item(15, 1).EntireRow.Delete()

# Find the new last row of data after the deletion.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

# After the shift, recompute the H (Numb) and I (RF) columns for the block
# that used to carry RF = 12.68054545454545: H values that were 0 become -1,
# and the RF value itself is refreshed to 12.68472727272727.
for ($r = 20; $r -le $lastRow; $r++) {
    $h = $ws.Cells.Item($r, 8).Value2
    if ($h -eq 0) {
        $ws.Cells.Item($r, 8).Value = -1
    }
    $ws.Cells.Item($r, 9).Value = 12.68472727272727
}
